# Daily attendance processing - 2026-01-29 18:58:50
# For each row in the "Recorded By" column (G), if the value is exactly
# "System, <something>" (a two-item list with "System" listed first),
# swap the order so the other entry is listed first and "System" last.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.StartsWith("System, ")) {
        $parts = $val.Split(",")
        if ($parts.Count -eq 2) {
            $first = $parts[0].Trim()
            $second = $parts[1].Trim()
            if ($first -eq "System") {
                $cell.Value = "$second, $first"
            }
        }
    }
}
